# 05/05 - push code
# Re-header the "User" (display) and "Note" (lookup) sheets from the old
# Vietnamese / ad-hoc labels to the new English labels, and add a new
# "Full name" / "fullName" column (I) that is looked up from Note.

$wb = $excel.ActiveWorkbook

$wsUser = $wb.Worksheets.Item("User")
$wsNote = $wb.Worksheets.Item("Note")

# ---------------------------------------------------------------------
# 1. "User" sheet - visible, friendly column headers (row 1)
# ---------------------------------------------------------------------
$wsUser.Range("A1").Value = "ID"
$wsUser.Range("B1").Value = "Username"
$wsUser.Range("C1").Value = "Address"
$wsUser.Range("D1").Value = "Gender"
$wsUser.Range("E1").Value = "Email"
$wsUser.Range("F1").Value = "Phone number"
$wsUser.Range("G1").Value = "Day of birth"
$wsUser.Range("H1").Value = "Job title"
$wsUser.Range("I1").Value = "Full name"

# Give the new column roughly the same width as its neighbours.
$wsUser.Columns.Item(9).ColumnWidth = 16.6

# New lookup column (I) for every data row, mirroring the existing
# B:H VLOOKUP formulas but pointing at the new 9th Note column.
for ($r = 2; $r -le 100; $r++) {
    $wsUser.Range("I$r").Formula = "=IF(A$r=`"`",`"`",VLOOKUP(A$r,Note!`$A`$1:`$I`$816,9,0))"
}

# ---------------------------------------------------------------------
# 2. "Note" sheet - raw field-name headers (row 1) used by the VLOOKUPs
# ---------------------------------------------------------------------
$wsNote.Range("A1").Value = "ID"
$wsNote.Range("B1").Value = "username"
$wsNote.Range("C1").Value = "sex"
$wsNote.Range("D1").Value = "address"
$wsNote.Range("E1").Value = "emailAddress"
$wsNote.Range("F1").Value = "phoneNumber"
$wsNote.Range("G1").Value = "dateOfBirth"
$wsNote.Range("H1").Value = "jobTitle"
$wsNote.Range("I1").Value = "fullName"

# ---------------------------------------------------------------------
# 3. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$wsUser.Range("F9").Select()
$wsNote.Activate()
$wsNote.Range("J11").Select()
